$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = 1000
$ws.Range("F11").Value = "CRCW06031K00JNEA"
$ws.Range("G11").Value = "https://octopart.com/crcw06031k00jnea-vishay-55388219?r=sp"
